# Update market price / profit data cells across sheets per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 72719.71000000001
$ws.Range("I21").Value = 73807.60000000001
$ws.Range("J21").Value = 70000
$ws.Range("K21").Value = 73807.60000000001
$ws.Range("L21").Value = 70000
$ws.Range("M21").Value = -73339.60000000001
$ws.Range("N21").Value = -70936

$ws.Range("H23").Value = 72719.71000000001
$ws.Range("I23").Value = 73807.60000000001
$ws.Range("J23").Value = 70000
$ws.Range("K23").Value = 73807.60000000001
$ws.Range("L23").Value = 70000
$ws.Range("M23").Value = -73573.60000000001
$ws.Range("N23").Value = -70468

$ws.Range("H29").Value = 1170
$ws.Range("I29").Value = 450
$ws.Range("J29").Value = 1350
$ws.Range("K29").Value = 1350
$ws.Range("L29").Value = 4050
$ws.Range("M29").Value = -1069
$ws.Range("N29").Value = -4612

$ws.Range("H43").Value = 794.7
$ws.Range("I43").Value = 625.25
$ws.Range("J43").Value = 837.0625
$ws.Range("K43").Value = 625.25
$ws.Range("L43").Value = 837.0625
$ws.Range("M43").Value = -556.25
$ws.Range("N43").Value = -975.0625

$ws.Range("H129").Value = 1851.0358
$ws.Range("I129").Value = 623.36365
$ws.Range("J129").Value = 2645.4119
$ws.Range("K129").Value = 1870.09095
$ws.Range("L129").Value = 7936.2357
$ws.Range("M129").Value = 3129.90905
$ws.Range("N129").Value = -17936.2357

$ws.Range("H132").Value = 4690.3677
$ws.Range("I132").Value = 3777.4314
$ws.Range("J132").Value = 7429.1763
$ws.Range("K132").Value = 11332.2942
$ws.Range("L132").Value = 22287.5289
$ws.Range("M132").Value = -8802.2942
$ws.Range("N132").Value = -27347.5289

$ws.Range("H138").Value = 2482.427
$ws.Range("I138").Value = 1633.0322
$ws.Range("J138").Value = 2887.5232
$ws.Range("K138").Value = 4899.096600000001
$ws.Range("L138").Value = 8662.569600000001
$ws.Range("M138").Value = 240.9033999999992
$ws.Range("N138").Value = -18942.5696

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3765
$ws.Range("I32").Value = 3765
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3765
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3478

$ws.Range("H61").Value = 1640.6545
$ws.Range("I61").Value = 1469.3103
$ws.Range("J61").Value = 1831.7693
$ws.Range("K61").Value = 1469.3103
$ws.Range("L61").Value = 1831.7693
$ws.Range("M61").Value = -1257.3103
$ws.Range("N61").Value = -2255.7693

$ws.Range("H132").Value = 3212545.5
$ws.Range("I132").Value = 9759.571
$ws.Range("J132").Value = 5006105.5
$ws.Range("K132").Value = 29278.713
$ws.Range("L132").Value = 15018316.5
$ws.Range("M132").Value = -26748.713
$ws.Range("N132").Value = -15023376.5

$ws.Range("H136").Value = 1640.6545
$ws.Range("I136").Value = 1469.3103
$ws.Range("J136").Value = 1831.7693
$ws.Range("K136").Value = 4407.9309
$ws.Range("L136").Value = 5495.3079
$ws.Range("M136").Value = -1857.9309
$ws.Range("N136").Value = -10595.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 14499.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 14499.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 14499.5
$ws.Range("N88").Value = -15311.5

$ws.Range("H91").Value = 14499.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 14499.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 14499.5
$ws.Range("N91").Value = -17307.5

$ws.Range("H96").Value = 14929.333
$ws.Range("I96").Value = 6952
$ws.Range("J96").Value = 22906.666
$ws.Range("K96").Value = 6952
$ws.Range("L96").Value = 22906.666
$ws.Range("M96").Value = -4206
$ws.Range("N96").Value = -28398.666

$ws.Range("H109").Value = 31600
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 31600
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 31600
$ws.Range("N109").Value = -34374

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2214.45
$ws.Range("I31").Value = 1163.0728
$ws.Range("J31").Value = 3499.4666
$ws.Range("K31").Value = 1163.0728
$ws.Range("L31").Value = 3499.4666
$ws.Range("M31").Value = -868.0727999999999
$ws.Range("N31").Value = -4089.4666

$ws.Range("H34").Value = 2214.45
$ws.Range("I34").Value = 1163.0728
$ws.Range("J34").Value = 3499.4666
$ws.Range("K34").Value = 1163.0728
$ws.Range("L34").Value = 3499.4666
$ws.Range("M34").Value = -961.0727999999999
$ws.Range("N34").Value = -3903.4666

$ws.Range("H43").Value = 31600
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 31600
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 31600
$ws.Range("N43").Value = -31968

$ws.Range("H101").Value = 31600
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 31600
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 31600
$ws.Range("N101").Value = -38090

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4388098
$ws.Range("I126").Value = 8335143.5
$ws.Range("J126").Value = 2492.4443
$ws.Range("K126").Value = 25005430.5
$ws.Range("L126").Value = 7477.3329
$ws.Range("M126").Value = -25002960.5
$ws.Range("N126").Value = -12417.3329

$ws.Range("H133").Value = 49556
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49556
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49556
$ws.Range("N133").Value = -59676

$ws.Range("H137").Value = 48000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 48000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 48000
$ws.Range("N137").Value = -58200

$ws.Range("H140").Value = 31666.666
$ws.Range("I140").Value = 31666.666
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 31666.666
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -26486.666
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 17523.092
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 17523.092
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 17523.092
$ws.Range("N97").Value = -19505.092

$ws.Range("H101").Value = 18246.1
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 18246.1
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 18246.1
$ws.Range("N101").Value = -24736.1

$ws.Range("H111").Value = 40387
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 40387
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 40387
$ws.Range("N111").Value = -48567

$ws.Range("H122").Value = 18309.715
$ws.Range("I122").Value = 27847
$ws.Range("J122").Value = 5593.3335
$ws.Range("K122").Value = 83541
$ws.Range("L122").Value = 16780.0005
$ws.Range("M122").Value = -81091
$ws.Range("N122").Value = -21680.0005

$ws.Range("H132").Value = 18217.514
$ws.Range("I132").Value = 20770.254
$ws.Range("J132").Value = 8176.7334
$ws.Range("K132").Value = 62310.762
$ws.Range("L132").Value = 24530.2002
$ws.Range("M132").Value = -59780.762
$ws.Range("N132").Value = -29590.2002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10626
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 10626
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 10626
$ws.Range("N45").Value = -11608

$ws.Range("H68").Value = 29249.75
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 29249.75
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 29249.75
$ws.Range("N68").Value = -30871.75

$ws.Range("H69").Value = 9085.1
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 9085.1
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 9085.1
$ws.Range("N69").Value = -10583.1

$ws.Range("H71").Value = 29249.75
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 29249.75
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 87749.25
$ws.Range("N71").Value = -95861.25

$ws.Range("H72").Value = 9085.1
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 9085.1
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 27255.3
$ws.Range("N72").Value = -34743.3

$ws.Range("H98").Value = 22600
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 22600
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 22600
$ws.Range("N98").Value = -28590

$ws.Range("H108").Value = 39500
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 39500
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 39500
$ws.Range("N108").Value = -47180

$ws.Range("H122").Value = 2643.65
$ws.Range("I122").Value = 2432.2727
$ws.Range("J122").Value = 2902
$ws.Range("K122").Value = 7296.8181
$ws.Range("L122").Value = 8706
$ws.Range("M122").Value = -4846.8181
$ws.Range("N122").Value = -13606

$ws.Range("H126").Value = 3954.7273
$ws.Range("I126").Value = 3389.111
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 10167.333
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -7697.332999999999
$ws.Range("N126").Value = -24440

$ws.Range("H132").Value = 1765.4
$ws.Range("I132").Value = 1362.4138
$ws.Range("J132").Value = 2214.8845
$ws.Range("K132").Value = 4087.2414
$ws.Range("L132").Value = 6644.6535
$ws.Range("M132").Value = -1557.2414
$ws.Range("N132").Value = -11704.6535
